{"js": "// \"Updated to 5th grade\" \u2014 the applicant is no longer a 4th-year EPF de\n// Cachan student applying for a Master's (4th-year) placement; the\n// cover letter is tweaked to drop the campus reference and to use the\n// (lower-case) \"master's\" wording throughout, plus a couple of small\n// wording/formatting touch-ups in the body paragraphs.\n\n// 1) \"As a fourth-year student at EPF de Cachan,\" -> \"...at EPF,\"\nlet results = context.document.body.search(\"EPF de Cachan,\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"EPF,\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) \"...any additional information in order to be a part of...\" ->\n//    \"...any additional information to be a part of...\"\nresults = context.document.body.search(\"information in order to be a part\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"information to be a part\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 3) \"Master's\" -> \"master's\" (lower-case) everywhere it appears. Only\n// the \"Master\" token is replaced (leaving the trailing \"'s\" as-is) so\n// the original straight apostrophe is preserved exactly.\nresults = context.document.body.search(\"Master\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"master\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 4) Justify the four main body paragraphs of the letter.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst targets = [\n  \"fourth-year student\",\n  \"My education in computer engineering\",\n  \"Indeed, I aspire\",\n  \"I sincerely appreciate your consideration\",\n];\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (targets.some((t) => text.includes(t))) {\n    paragraphs.items[i].alignment = Word.Alignment.justified;\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Text edits (use narrow Find/Replace scopes so we never touch the\n# straight apostrophe in \"Master's\"/\"master's\" \u2014 replacing the whole\n# word including the apostrophe causes AutoCorrect to turn it into a\n# curly quote). ---\n\n# \"EPF de Cachan,\" -> \"EPF,\" (student no longer attends the Cachan campus)\n$find = $d.Content.Find\n$find.Text = \"EPF de Cachan,\"\n$find.Replacement.Text = \"EPF,\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# \"in order to be a part\" -> \"to be a part\"\n$find = $d.Content.Find\n$find.Text = \"information in order to be a part\"\n$find.Replacement.Text = \"information to be a part\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# \"Master's\" -> \"master's\" (lowercase) in all three occurrences. Replace\n# only the \"Master\" token so the trailing \"'s\" (straight apostrophe) is\n# left completely untouched.\n$find = $d.Content.Find\n$find.Text = \"Master\"\n$find.Replacement.Text = \"master\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# --- Paragraph formatting: justify the body paragraphs of the letter ---\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*fourth-year student*\" -or `\n        $t -like \"*My education in computer engineering*\" -or `\n        $t -like \"*Indeed, I aspire*\" -or `\n        $t -like \"*I sincerely appreciate*\") {\n        $p.Alignment = 3\n    }\n}\n"}
